$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.910.71'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.48%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.908.17'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.37%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.11'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.13%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4585'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.84%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.64%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07717'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9798'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.06'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.39%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.918.27'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.84%  '

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.665'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.80%  '

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.937'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.05%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07043'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.17%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '83.76'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.32%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009451'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.64%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.64%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '28.898.68'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.322'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.01%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.51%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.16%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.56'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.51%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.99'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.16%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.80%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '117.79'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.57%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.866'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.40%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09274'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.74%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.8649'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.71%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.075'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.82%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.45%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.089'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05714'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.95%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.164'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.84%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02043'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.65%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.56%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5484'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.24%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1752'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.49%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.886'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.58%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.311'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.19%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5157'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.32%  '

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.106'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.63%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '11.20'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.76%  '

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.06890'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.777'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.23%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '110.39'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.81%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000002545'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -13.91%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.2863'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.34%  '
